$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 21-49 (excluding the special H-change rows) just get a new RF value in column I.
$rfValue = 66.00500000000001
$plainRows = @(21,22,23,24,25,26,27,29,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49)
foreach ($r in $plainRows) {
    $ws.Cells.Item($r, 9).Value = $rfValue
}

# Rows 28, 30, 50, 51 also change H from 0 to -1, besides adding the RF value in column I.
$changedRows = @(28,30,50,51)
foreach ($r in $changedRows) {
    $ws.Cells.Item($r, 8).Value = -1
    $ws.Cells.Item($r, 9).Value = $rfValue
}

# Row 66 ("Scyliorhinus canicula") is removed entirely; subsequent rows shift up,
# and the dimension shrinks from A1:K71 to A1:K70.
$ws.Rows("66").Delete()
